{"js": "// Update the stack-trace line numbers embedded in the document text so\n// they match the new M2DocEvaluator.java / M2DocUtils.java /\n// AbstractTemplatesTestSuite.java source positions (the custom docProp\n// \"m:var:self\" version bump shifted every later line by a few lines).\nconst replacements = [\n  [\"M2DocEvaluator.java:1003\", \"M2DocEvaluator.java:1049\"],\n  [\"M2DocEvaluator.java:1038\", \"M2DocEvaluator.java:1084\"],\n  [\"M2DocEvaluator.java:1254\", \"M2DocEvaluator.java:1300\"],\n  [\"M2DocEvaluator.java:275)\", \"M2DocEvaluator.java:278)\"],\n  [\"M2DocEvaluator.java:264)\", \"M2DocEvaluator.java:267)\"],\n  [\"M2DocUtils.java:712\", \"M2DocUtils.java:694\"],\n  [\"AbstractTemplatesTestSuite.java:459\", \"AbstractTemplatesTestSuite.java:475\"],\n  [\"AbstractTemplatesTestSuite.java:369\", \"AbstractTemplatesTestSuite.java:384\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the stack-trace line numbers embedded in the document text so\n# they match the new M2DocEvaluator.java / M2DocUtils.java /\n# AbstractTemplatesTestSuite.java source positions (the custom docProp\n# \"m:var:self\" version bump shifted every later line by a few lines).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"M2DocEvaluator.java:1003\", \"M2DocEvaluator.java:1049\"),\n    @(\"M2DocEvaluator.java:1038\", \"M2DocEvaluator.java:1084\"),\n    @(\"M2DocEvaluator.java:1254\", \"M2DocEvaluator.java:1300\"),\n    @(\"M2DocEvaluator.java:275)\", \"M2DocEvaluator.java:278)\"),\n    @(\"M2DocEvaluator.java:264)\", \"M2DocEvaluator.java:267)\"),\n    @(\"M2DocUtils.java:712\", \"M2DocUtils.java:694\"),\n    @(\"AbstractTemplatesTestSuite.java:459\", \"AbstractTemplatesTestSuite.java:475\"),\n    @(\"AbstractTemplatesTestSuite.java:369\", \"AbstractTemplatesTestSuite.java:384\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
